$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.6226591760299626
$summary.Range("C2").Value = 0.5746864310148233
$summary.Range("D2").Value = 0.9438202247191011
$summary.Range("E2").Value = 0.7143869596031184
$summary.Range("F2").Value = 0.8363757052771325
$summary.Range("G2").Value = 0.921065579531876
$summary.Range("H2").Value = 0.7968410273674761
$summary.Range("I2").Value = 504
$summary.Range("J2").Value = 373
$summary.Range("K2").Value = 161
$summary.Range("L2").Value = 30

# ---------------------------------------------------------------
# Sheet 2: Classification Report
# ---------------------------------------------------------------
$classRep = $wb.Worksheets.Item("Classification Report")

$classRep.Range("B2").Value = 0.8429319371727748
$classRep.Range("C2").Value = 0.301498127340824
$classRep.Range("D2").Value = 0.4441379310344827

$classRep.Range("B3").Value = 0.5746864310148233
$classRep.Range("C3").Value = 0.9438202247191011
$classRep.Range("D3").Value = 0.7143869596031184

$classRep.Range("B4").Value = 0.6226591760299626
$classRep.Range("C4").Value = 0.6226591760299626
$classRep.Range("D4").Value = 0.6226591760299626
$classRep.Range("E4").Value = 0.6226591760299626

$classRep.Range("B5").Value = 0.7088091840937991
$classRep.Range("C5").Value = 0.6226591760299626
$classRep.Range("D5").Value = 0.5792624453188006

$classRep.Range("B6").Value = 0.7088091840937991
$classRep.Range("C6").Value = 0.6226591760299626
$classRep.Range("D6").Value = 0.5792624453188006

# ---------------------------------------------------------------
# Sheet 3: Confusion Matrix
# ---------------------------------------------------------------
$confMatrix = $wb.Worksheets.Item("Confusion Matrix")

$confMatrix.Range("B2").Value = 161
$confMatrix.Range("C2").Value = 373

$confMatrix.Range("B3").Value = 30
$confMatrix.Range("C3").Value = 504
